$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4002
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4004
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 4004
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4500
$ws.Range("H67").Value = 4002
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4004
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 4004
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5720
$ws.Range("H135").Value = 879.2222
$ws.Range("I135").Value = 872.1177
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 7849.0593
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -5314.0593
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 3707425.5
$ws.Range("J137").Value = 2671
$ws.Range("L137").Value = 8013
$ws.Range("N137").Value = -13113

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3263.31
$ws.Range("I32").Value = 2775.353
$ws.Range("J32").Value = 6028.4
$ws.Range("K32").Value = 2775.353
$ws.Range("L32").Value = 6028.4
$ws.Range("M32").Value = -2488.353
$ws.Range("N32").Value = -6602.4
$ws.Range("H61").Value = 1976.5
$ws.Range("I61").Value = 1149.9166
$ws.Range("J61").Value = 2968.4
$ws.Range("K61").Value = 1149.9166
$ws.Range("L61").Value = 2968.4
$ws.Range("M61").Value = -937.9166
$ws.Range("N61").Value = -3392.4
$ws.Range("H74").Value = 637.75
$ws.Range("I74").Value = 563.9
$ws.Range("J74").Value = 1007
$ws.Range("K74").Value = 563.9
$ws.Range("L74").Value = 1007
$ws.Range("M74").Value = 310.1
$ws.Range("N74").Value = -2755
$ws.Range("H77").Value = 637.75
$ws.Range("I77").Value = 563.9
$ws.Range("J77").Value = 1007
$ws.Range("K77").Value = 2819.5
$ws.Range("L77").Value = 5035
$ws.Range("M77").Value = 1548.5
$ws.Range("N77").Value = -13771
$ws.Range("H132").Value = 33337396
$ws.Range("I132").Value = 47622468
$ws.Range("J132").Value = 5558.4443
$ws.Range("K132").Value = 142867404
$ws.Range("L132").Value = 16675.3329
$ws.Range("M132").Value = -142864874
$ws.Range("N132").Value = -21735.3329
$ws.Range("H136").Value = 1976.5
$ws.Range("I136").Value = 1149.9166
$ws.Range("J136").Value = 2968.4
$ws.Range("K136").Value = 3449.7498
$ws.Range("L136").Value = 8905.200000000001
$ws.Range("M136").Value = -899.7498000000001
$ws.Range("N136").Value = -14005.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 33295
$ws.Range("J69").Value = 33295
$ws.Range("L69").Value = 33295
$ws.Range("N69").Value = -34917
$ws.Range("H72").Value = 33295
$ws.Range("J72").Value = 33295
$ws.Range("L72").Value = 99885
$ws.Range("N72").Value = -107997
$ws.Range("H134").Value = 4409.391
$ws.Range("I134").Value = 4319.1763
$ws.Range("K134").Value = 12957.5289
$ws.Range("M134").Value = -10422.5289

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2130161.2
$ws.Range("I31").Value = 2779655
$ws.Range("J31").Value = 4545.4546
$ws.Range("K31").Value = 2779655
$ws.Range("L31").Value = 4545.4546
$ws.Range("M31").Value = -2779360
$ws.Range("N31").Value = -5135.4546
$ws.Range("H34").Value = 2130161.2
$ws.Range("I34").Value = 2779655
$ws.Range("J34").Value = 4545.4546
$ws.Range("K34").Value = 2779655
$ws.Range("L34").Value = 4545.4546
$ws.Range("M34").Value = -2779453
$ws.Range("N34").Value = -4949.4546
$ws.Range("H132").Value = 3483.8667
$ws.Range("I132").Value = 2688.1667
$ws.Range("K132").Value = 8064.500100000001
$ws.Range("M132").Value = -5534.500100000001
$ws.Range("H134").Value = 1932.4117
$ws.Range("I134").Value = 990.06665
$ws.Range("K134").Value = 2970.19995
$ws.Range("M134").Value = -435.1999500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1286.7693
$ws.Range("I5").Value = 580.8889
$ws.Range("J5").Value = 2875
$ws.Range("K5").Value = 1742.6667
$ws.Range("L5").Value = 8625
$ws.Range("M5").Value = -1630.6667
$ws.Range("N5").Value = -8849
$ws.Range("H23").Value = 170
$ws.Range("J23").Value = 192.5
$ws.Range("L23").Value = 577.5
$ws.Range("N23").Value = -1047.5
$ws.Range("H122").Value = 1001.6957
$ws.Range("I122").Value = 359.23077
$ws.Range("J122").Value = 1836.9
$ws.Range("K122").Value = 3233.07693
$ws.Range("L122").Value = 16532.1
$ws.Range("M122").Value = -783.0769300000002
$ws.Range("N122").Value = -21432.1
$ws.Range("H132").Value = 3842
$ws.Range("J132").Value = 4196.75
$ws.Range("L132").Value = 37770.75
$ws.Range("N132").Value = -42830.75
$ws.Range("H135").Value = 1286.7693
$ws.Range("I135").Value = 580.8889
$ws.Range("J135").Value = 2875
$ws.Range("K135").Value = 5228.0001
$ws.Range("L135").Value = 25875
$ws.Range("M135").Value = -2693.0001
$ws.Range("N135").Value = -30945
$ws.Range("H139").Value = 5765.032
$ws.Range("I139").Value = 2141.3
$ws.Range("J139").Value = 12353.637
$ws.Range("K139").Value = 6423.900000000001
$ws.Range("L139").Value = 37060.911
$ws.Range("M139").Value = -1283.900000000001
$ws.Range("N139").Value = -47340.911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2937.7
$ws.Range("I132").Value = 2288.8462
$ws.Range("J132").Value = 4142.7144
$ws.Range("K132").Value = 6866.5386
$ws.Range("L132").Value = 12428.1432
$ws.Range("M132").Value = -4336.5386
$ws.Range("N132").Value = -17488.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3622.3044
$ws.Range("I132").Value = 1963.75
$ws.Range("J132").Value = 4506.8667
$ws.Range("K132").Value = 5891.25
$ws.Range("L132").Value = 13520.6001
$ws.Range("M132").Value = -3361.25
$ws.Range("N132").Value = -18580.6001
$ws.Range("H136").Value = 2633848.2
$ws.Range("I136").Value = 3450128.8
$ws.Range("K136").Value = 10350386.4
$ws.Range("M136").Value = -10347836.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 262253.5
$ws.Range("I132").Value = 360013.84
$ws.Range("J132").Value = 13408.909
$ws.Range("K132").Value = 1080041.52
$ws.Range("L132").Value = 40226.727
$ws.Range("M132").Value = -1077511.52
$ws.Range("N132").Value = -45286.727
$ws.Range("H136").Value = 1638.9
$ws.Range("I136").Value = 948.4286
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 2845.2858
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -295.2857999999997
$ws.Range("N136").Value = -14850

Write-Output "done"
